$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (unchanged values, kept for completeness)
$ws.Range("A1").Value = "Variable/Konfiguration"
$ws.Range("B1").Value = "Config1"
$ws.Range("C1").Value = "Config2"

# Row 2
$ws.Range("A2").Value = "BaseStateChromium"
$ws.Range("B2").Value = "Started"
$ws.Range("C2").Value = "Not Started"

# Row 3
$ws.Range("A3").Value = "BaseStatePixel9Pro_API35"
$ws.Range("B3").Value = "Not Started"
$ws.Range("C3").Value = "Started"

# New row 4
$ws.Range("A4").Value = "AUT"
$ws.Range("B4").Value = "Chromium"
$ws.Range("C4").Value = "Pixel9Pro_API35"

# Adjust column widths to match autofit ("bestFit") sizing seen in target
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Update selection to reflect the new active cell location from the diff
$ws.Range("P18").Select() | Out-Null
